$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '301.32'
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '-0.63%'
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '36.47'
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '2.55%'
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '4.981'
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '-2.16%'
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.07734'
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '-0.11%'
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '2.070'
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '-6.78%'
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '7.903'
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '-1.62%'
$cell.Style = "Normal"

$cell = $ws.Range("B8")
$cell.NumberFormat = "@"
$cell.Value = 'MXToken'
$cell.Style = "Normal"

$cell = $ws.Range("C8")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.9221'
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '-0.57%'
$cell.Style = "Normal"

$cell = $ws.Range("B9")
$cell.NumberFormat = "@"
$cell.Value = 'LiechtensteinCryptoassetsExchange'
$cell.Style = "Normal"

$cell = $ws.Range("C9")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.09733'
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '1.44%'
$cell.Style = "Normal"

$cell = $ws.Range("B10")
$cell.NumberFormat = "@"
$cell.Value = 'WazirX'
$cell.Style = "Normal"

$cell = $ws.Range("C10")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.1852'
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '1.07%'
$cell.Style = "Normal"

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = 'MandalaExchangeToken'
$cell.Style = "Normal"

$cell = $ws.Range("C11")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.08565'
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '0.11%'
$cell.Style = "Normal"

$cell = $ws.Range("B12")
$cell.NumberFormat = "@"
$cell.Value = 'BitrueCoin'
$cell.Style = "Normal"

$cell = $ws.Range("C12")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.03509'
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '-3.69%'
$cell.Style = "Normal"

$cell = $ws.Range("B13")
$cell.NumberFormat = "@"
$cell.Value = 'BitMartToken'
$cell.Style = "Normal"

$cell = $ws.Range("C13")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.09941'
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '-0.30%'
$cell.Style = "Normal"

$cell = $ws.Range("B14")
$cell.NumberFormat = "@"
$cell.Value = 'BitForexToken'
$cell.Style = "Normal"

$cell = $ws.Range("C14")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.001465'
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '-0.97%'
$cell.Style = "Normal"

$cell = $ws.Range("B15")
$cell.NumberFormat = "@"
$cell.Value = 'TigerCash'
$cell.Style = "Normal"

$cell = $ws.Range("C15")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.005631'
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '-1.17%'
$cell.Style = "Normal"

$cell = $ws.Range("B16")
$cell.NumberFormat = "@"
$cell.Value = 'LEO'
$cell.Style = "Normal"

$cell = $ws.Range("C16")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.467'
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '-0.32%'
$cell.Style = "Normal"

$cell = $ws.Range("B17")
$cell.NumberFormat = "@"
$cell.Value = 'GateToken'
$cell.Style = "Normal"

$cell = $ws.Range("C17")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '4.023'
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '-0.26%'
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '2.281'
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '4.40%'
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.3408'
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '-1.59%'
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.1341'
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '1.25%'
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '4.767'
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '4.38%'
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.2195'
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '-2.12%'
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.04595'
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.005086'
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '12.88%'
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.001229'
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '-0.63%'
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.0001398'
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '6.87%'
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.01763'
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '0.30%'
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.04641'
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '-1.65%'
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.007435'
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '-6.49%'
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.1390'
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '-1.26%'
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.007700'
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '0.32%'
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.002246'
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '0.86%'
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.01032'
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '6.81%'
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00006168'
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '-1.17%'
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000748'
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '-0.90%'
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0005785'
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '-0.26%'
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '35.68'
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '516.61%'
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.001995'
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '-26.19%'
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.00002094'
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '-0.90%'
$cell.Style = "Normal"
